$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Pelanggan test case): kodebarang=GP01, netto="ABC", satuanberat=<blank>, harga=1001
$ws.Range("B3").Value = "GP01"
$ws.Range("C3").Value = "ABC"
$ws.Range("D3").Borders.LineStyle = -4142
$ws.Range("E3").Value = 1001

# Row 4 (Salesman test case): namabarang=Gula Pasir, satuanberat=Kilogram, harga=TEST1002
$ws.Range("A4").Value = "Gula Pasir"
$ws.Range("D4").Value = "Kilogram"
$ws.Range("E4").Value = "TEST1002"

# Row 5 (Supplier test case): namabarang=Gula Aren, kodebarang=GP01, netto=Asd1fgh, satuanberat=Kilogram, harga=TEST
$ws.Range("A5").Value = "Gula Aren"
$ws.Range("B5").Value = "GP01"
$ws.Range("C5").Value = "Asd1fgh"
$ws.Range("D5").Value = "Kilogram"
$ws.Range("E5").Value = "TEST"

[void]$ws.Range("E5").Select()
